# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (column F) values for the affected rows to reflect the repulled data
$ws.Range("F4").Value  = -1
$ws.Range("F10").Value = -1
$ws.Range("F13").Value = -4
$ws.Range("F16").Value = -6
$ws.Range("F18").Value = -2
$ws.Range("F21").Value = -3
$ws.Range("F23").Value = -7
$ws.Range("F24").Value = -6
$ws.Range("F28").Value = 7
$ws.Range("F33").Value = -4
